# chore: add shifts 12/14 - 12/24
#
# Appends 7 new daily-shift rows (rows 221-227) below the existing data,
# continuing the Date / Cash Tips / Hours / Wage / Total Income table.
# Column E keeps the same "B+(C*D)" formula pattern as every prior row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shift data: Date (Excel serial), Cash Tips, Hours, Wage
$shifts = @(
    @(221, 44911, 260, 8.5,  15),
    @(222, 44912, 200, 6.5,  15),
    @(223, 44913, 613, 12.5, 15),
    @(224, 44915, 105, 7,    15),
    @(225, 44916, 230, 7,    15),
    @(226, 44918, 235, 8,    15),
    @(227, 44919, 125, 3,    15)
)

# Row 220 is the last existing data row - carry its formatting (styles +
# row height) down into each freshly-appended row before filling values.
foreach ($shift in $shifts) {
    $row = $shift[0]

    $ws.Range("A220:E220").Copy()
    $ws.Range("A" + $row + ":E" + $row).PasteSpecial(-4122)
    $ws.Rows.Item($row).RowHeight = 13.55

    $ws.Cells.Item($row, 1).Value = $shift[1]
    $ws.Cells.Item($row, 2).Value = $shift[2]
    $ws.Cells.Item($row, 3).Value = $shift[3]
    $ws.Cells.Item($row, 4).Value = $shift[4]
    $ws.Cells.Item($row, 5).Formula = "=B" + $row + "+(C" + $row + "*D" + $row + ")"
}

$excel.CutCopyMode = 0
